# Update the marksheet's "Correct / Total" marks.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("quiz")

# Right answers count (row 11, "Marking" row, column B)
$ws.Range("B11").Value = 5

# Total score (row 12, "Total" row, column B) and the "Corr/total" display text (E12)
$ws.Range("B12").Value = 120
$ws.Range("E12").Value = "120/140"
